$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Cells.Item(20, 1).Value = 19
    $ws.Cells.Item(20, 2).NumberFormat = "@"
    $ws.Cells.Item(20, 2).Value = "2026-02-16"
    $ws.Cells.Item(20, 3).NumberFormat = "@"
    $ws.Cells.Item(20, 3).Value = "22:59:15"
    $ws.Cells.Item(20, 4).Value = "base_strategy"
    $ws.Cells.Item(20, 5).Value = "UP"
    $ws.Cells.Item(20, 6).Value = 0.5
    $ws.Cells.Item(20, 7).Value = ""
    $ws.Cells.Item(20, 8).Value = "OPEN"
    $ws.Cells.Item(20, 9).Value = 0
    $ws.Cells.Item(20, 10).Value = 0
    $ws.Cells.Item(20, 11).Value = 100
    $ws.Cells.Item(20, 12).Value = 0
    $ws.Cells.Item(20, 13).Value = 0
    $ws.Cells.Item(20, 14).Value = 0.6
    $ws.Cells.Item(20, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(20, 16).Value = ""
    $ws.Cells.Item(20, 17).Value = 0
}
